$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 22 (1-indexed) is the "Should we meet with the HET (diabetes) team..." row.
# We need to keep that row but empty out all 4 cells (reduce each to a single,
# completely bare paragraph), and delete the 3 following rows entirely
# ("Should we scale the graph...", "Should we have a new logo...",
#  "Should we show the sensors...").

$hetRow = $t.Rows.Item(22)

# Insert a brand-new, genuinely blank row immediately before the HET row.
# Newly added rows have simple empty paragraphs (no leftover formatting),
# which is exactly what the target XML wants.
$blankRow = $t.Rows.Add($hetRow)

# Give each cell of the blank row some throwaway text and then delete it via an
# absolute document Range; this reliably collapses the cell back down to a
# single bare <w:p/> in this runtime (directly deleting an already-empty
# placeholder run is a no-op here).
for ($c = 1; $c -le 4; $c++) {
    $cell = $blankRow.Cells.Item($c)
    $cell.Range.Text = "x"
    $r = $cell.Range
    $d.Range($r.Start, $r.End).Delete()
}

# The original HET row (with its real content) is now pushed down to index 23,
# immediately after our new blank row. Delete it along with the next three
# rows (scale graph / new logo / sensors not working), which are now at
# indices 23-25 (after the original HET row at 23 was removed once, the
# following rows shift up, so we just keep deleting row 23 four times).
for ($i = 0; $i -lt 4; $i++) {
    $t.Rows.Item(23).Delete()
}
